$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new employee ("michael") was appended below the existing (sorted)
# table, in the same "role"/count layout as everybody else.
$ws.Range("A26").Value = "michael"
$ws.Range("B26").Value = "מלקט"
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0

# Match the formatting of the row above (borders/alignment/number format)
# so the new row looks like every other data row in the table.
$ws.Range("A25:F25").Copy()
$ws.Range("A26:F26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Leave the selection where the user ended up after adding the row.
$ws.Range("G29").Select()
